# Apply cryptos list update (Sat Aug 24 20:43:18 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.155.80"
$ws.Range("E2").Value = "  +1.50%  "

$ws.Range("D3").Value = "'2.789.32"
$ws.Range("E3").Value = "  +2.56%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'591.27"
$ws.Range("E5").Value = "  +0.68%  "

$ws.Range("D6").Value = "'160.80"
$ws.Range("E6").Value = "  +6.60%  "

$ws.Range("E7").Value = "  +2.55%  "

$ws.Range("E8").Value = "  +0.26%  "

$ws.Range("D9").Value = "'6.77"
$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("E10").Value = "  +1.67%  "

$ws.Range("E11").Value = "  +2.24%  "

$ws.Range("D12").Value = "'0.160"
$ws.Range("E12").Value = "  +0.99%  "

$ws.Range("D13").Value = "'3.281.67"
$ws.Range("E13").Value = "  +2.39%  "

$ws.Range("D14").Value = "'27.39"
$ws.Range("E14").Value = "  +2.38%  "

$ws.Range("D15").Value = "'64.067.67"
$ws.Range("E15").Value = "  +1.51%  "

$ws.Range("E16").Value = "  +6.18%  "

$ws.Range("D17").Value = "'2.798.76"
$ws.Range("E17").Value = "  +2.03%  "

$ws.Range("D18").Value = "'12.50"
$ws.Range("E18").Value = "  +4.19%  "

$ws.Range("D19").Value = "'5.04"
$ws.Range("E19").Value = "  +3.63%  "

$ws.Range("D20").Value = "'367.90"
$ws.Range("E20").Value = "  +1.10%  "

$ws.Range("D21").Value = "'7.04"

$ws.Range("E22").Value = "  +7.69%  "

$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").Value = "'67.34"
$ws.Range("E24").Value = "  +2.65%  "

$ws.Range("D25").Value = "'0.177"
$ws.Range("E25").Value = "  +6.33%  "

$ws.Range("E26").Value = "  +3.21%  "

$ws.Range("D27").Value = "'0.0₃0970"
$ws.Range("E27").Value = "  +12.28%  "

$ws.Range("E28").Value = "  +0.36%  "

$ws.Range("E29").Value = "  +1.96%  "

$ws.Range("E30").Value = "  +1.80%  "

$ws.Range("E31").Value = "  +6.74%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "'170.85"
$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'5.18"
$ws.Range("E33").Value = "  +8.74%  "

$ws.Range("D34").Value = "'20.93"
$ws.Range("E34").Value = "  +1.78%  "

$ws.Range("D36").Value = "'1.50"
$ws.Range("E36").Value = "  +4.10%  "

$ws.Range("E37").Value = "  +2.74%  "

$ws.Range("E38").Value = "  +2.10%  "

$ws.Range("D39").Value = "'343.79"
$ws.Range("E39").Value = "  -2.26%  "

$ws.Range("E40").Value = "  +0.54%  "

$ws.Range("D41").Value = "'6.30"
$ws.Range("E41").Value = "  +10.65%  "

$ws.Range("D42").Value = "'40.33"
$ws.Range("E42").Value = "  +2.76%  "

$ws.Range("D43").Value = "'22.62"
$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("D44").Value = "'0.0614"
$ws.Range("E44").Value = "  +3.81%  "

$ws.Range("D45").Value = "'22.50"
$ws.Range("E45").Value = "  +3.04%  "

$ws.Range("D46").Value = "'0.654"
$ws.Range("E46").Value = "  +1.94%  "

$ws.Range("E47").Value = "  +1.69%  "

$ws.Range("D48").Value = "'138.57"
$ws.Range("E48").Value = "  -1.62%  "

$ws.Range("E49").Value = "  +2.33%  "

$ws.Range("D50").Value = "'2.178.67"
$ws.Range("E50").Value = "  +1.02%  "

$ws.Range("E51").Value = "  +0.36%  "
